# The workbook's single sheet contains a long list of daily "Mango" price
# records for "Vega Modelo de Temuco", one per row, running from row 2
# through row 618 (dimension A1:T618).
#
# This edit inserts a brand-new daily record as a new row 497, which
# pushes the previously-existing rows 497-618 down to rows 498-619
# (dimension becomes A1:T619). The new row carries the same static
# columns (Mercado/Región/Codreg/Tipo/Producto/Categoría/Variedad/
# Calidad/Unidad/Kg per unidad) as all of the other Mango rows, with its
# own date, volume, price and origin values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 497; everything at/after 497 shifts
# down by one row (this also grows the sheet dimension to A1:T619).
$ws.Rows.Item(497).Insert()

# Populate the new row 497 with the new record's data.
$ws.Range('A497').Value = 10
$ws.Range('B497').Value = 'Vega Modelo de Temuco'
$ws.Range('C497').Value = 'La Araucanía'
$ws.Range('D497').Value = 45135
$ws.Range('E497').Value = 9
$ws.Range('F497').Value = 'Fruta'
$ws.Range('G497').Value = 100108
$ws.Range('H497').Value = 'Tropicales y subtropicales'
$ws.Range('I497').Value = 100108002
$ws.Range('J497').Value = 'Mango'
$ws.Range('K497').Value = 'Sin especificar'
$ws.Range('L497').Value = 'Primera'
$ws.Range('M497').Value = 600
$ws.Range('N497').Value = 8000
$ws.Range('O497').Value = 8000
$ws.Range('P497').Value = 8000
$ws.Range('Q497').Value = '$/bandeja 4 kilos'
$ws.Range('R497').Value = 'Brasil'
$ws.Range('S497').Value = 2000
$ws.Range('T497').Value = 4
